# GT-N8013 support is finally here.
# Insert a new bold change-log entry right after the
# "October 5, 2019 ovation: ... reboot." paragraph, matching that
# paragraph's formatting (pStyle normal0, left indent 720, bold sz28 run).

$d = $word.ActiveDocument

$anchorText = "October 5, 2019 ovation: Fixed the SystemUI crashing      repeatedly after a reboot."
$newEntryText = "October 26, 2019 n8013: First build."

# Locate the anchor paragraph's index by scanning the Paragraphs collection
# (more robust than relying on a hard-coded paragraph number).
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*$anchorText*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    Write-Host "ERROR: anchor paragraph not found"
} else {
    # Insert a brand new paragraph right after the anchor paragraph; Word
    # automatically carries over the anchor's paragraph/run formatting.
    $anchorRange = $d.Paragraphs.Item($anchorIndex).Range
    $anchorRange.Collapse(0)
    $anchorRange.InsertParagraphAfter()

    # Fill in the text of the newly created paragraph.
    $newParagraph = $d.Paragraphs.Item($anchorIndex + 1)
    $newParagraph.Range.Text = $newEntryText

    Write-Host "Inserted '$newEntryText' after paragraph $anchorIndex"
}
